$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shared strings (sharedStrings.xml) for Volume/date header text ---
$ws2 = $wb.Worksheets.Item(1)

# --- Step 1: set NumberFormat to Text for cells becoming text placeholders ---
$ws.Range("C20").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"

# --- Step 2: assign values (numbers or text) ---
$ws.Range("D14").Value = "0"
$ws.Range("E14").Value = "***.*"
$ws.Range("D15").Value = "0"
$ws.Range("E15").Value = "***.*"
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -40
$ws.Range("I15").Value = 3
$ws.Range("K15").Value = -50
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -40
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -14.285714285714
$ws.Range("F16").Value = 21
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = 23.529411764705
$ws.Range("I16").Value = 29
$ws.Range("J16").Value = 22
$ws.Range("K16").Value = 31.818181818181
$ws.Range("L16").Value = -29.268292682926
$ws.Range("M16").Value = -9.375
$ws.Range("N16").Value = -79.861111111111
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = -31.034482758620
$ws.Range("I17").Value = 32
$ws.Range("J17").Value = 53
$ws.Range("K17").Value = -39.622641509434
$ws.Range("L17").Value = -31.914893617021
$ws.Range("M17").Value = 23.076923076923
$ws.Range("N17").Value = -17.948717948717
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -57.142857142857
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = -76.470588235294
$ws.Range("I18").Value = 5
$ws.Range("J18").Value = 18
$ws.Range("K18").Value = -72.222222222222
$ws.Range("L18").Value = -80.769230769230
$ws.Range("M18").Value = -84.375
$ws.Range("N18").Value = -98.366013071895
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = -30
$ws.Range("G19").Value = 59
$ws.Range("H19").Value = -27.118644067796
$ws.Range("I19").Value = 64
$ws.Range("J19").Value = 79
$ws.Range("K19").Value = -18.987341772151
$ws.Range("L19").Value = -49.206349206349
$ws.Range("M19").Value = 20.754716981132
$ws.Range("N19").Value = -48.8
$ws.Range("C20").Value = "0"
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = -28.571428571428
$ws.Range("J20").Value = 18
$ws.Range("K20").Value = -16.666666666666
$ws.Range("L20").Value = -51.612903225806
$ws.Range("M20").Value = -40
$ws.Range("N20").Value = -93.589743589743
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 46
$ws.Range("E21").Value = -30.434782608695
$ws.Range("F21").Value = 101
$ws.Range("G21").Value = 143
$ws.Range("H21").Value = -29.370629370629
$ws.Range("I21").Value = 148
$ws.Range("J21").Value = 198
$ws.Range("K21").Value = -25.252525252525
$ws.Range("L21").Value = -45.985401459854
$ws.Range("M21").Value = -13.450292397660
$ws.Range("N21").Value = -82.669789227166
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 300
$ws.Range("F22").Value = 6
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 6
$ws.Range("J22").Value = 6
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = -33.333333333333
$ws.Range("M22").Value = 20
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = -25
$ws.Range("F24").Value = 76
$ws.Range("G24").Value = 116
$ws.Range("H24").Value = -34.482758620689
$ws.Range("I24").Value = 112
$ws.Range("J24").Value = 147
$ws.Range("K24").Value = -23.809523809523
$ws.Range("L24").Value = -56.923076923076
$ws.Range("M24").Value = 0.900900900900
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -85.714285714285
$ws.Range("F25").Value = 19
$ws.Range("G25").Value = 62
$ws.Range("H25").Value = -69.354838709677
$ws.Range("I25").Value = 22
$ws.Range("J25").Value = 72
$ws.Range("K25").Value = -69.444444444444
$ws.Range("L25").Value = -86.335403726708
$ws.Range("C26").Value = 25
$ws.Range("D26").Value = 22
$ws.Range("E26").Value = 13.636363636363
$ws.Range("F26").Value = 62
$ws.Range("G26").Value = 82
$ws.Range("H26").Value = -24.390243902439
$ws.Range("I26").Value = 94
$ws.Range("J26").Value = 110
$ws.Range("K26").Value = -14.545454545454
$ws.Range("L26").Value = -10.476190476190
$ws.Range("M26").Value = 0
$ws.Range("D27").Value = "0"
$ws.Range("E27").Value = "***.*"
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 20
$ws.Range("I27").Value = 6
$ws.Range("K27").Value = -14.285714285714
$ws.Range("L27").Value = -25
$ws.Range("C28").Value = 4
$ws.Range("E28").Value = 300
$ws.Range("F28").Value = 6
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 9
$ws.Range("J28").Value = 7
$ws.Range("K28").Value = 28.571428571428
$ws.Range("L28").Value = -40
$ws.Range("M29").Value = -100
$ws.Range("M30").Value = -100
$ws.Range("C31").Value = 1
$ws.Range("F31").Value = 1
$ws.Range("I31").Value = 1
$ws.Range("L31").Value = -66.666666666666
$ws.Range("L33").Value = -50

# --- Step 3: fix styles by pasting formats from stable anchor cells ---
$ws.Range("A14").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("G14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("G14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("H14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("H14").Copy()
$ws.Range("M29").PasteSpecial(-4122)
$ws.Range("H14").Copy()
$ws.Range("M30").PasteSpecial(-4122)
$ws.Range("G14").Copy()
$ws.Range("C31").PasteSpecial(-4122)
$ws.Range("G14").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$ws.Range("G14").Copy()
$ws.Range("I31").PasteSpecial(-4122)

$excel.CutCopyMode = 0
